$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) and one data row (row 2):
#   A: "Unnamed: 0" (blank per-row), B: course code, C: material path
# Append a new data row (row 3) for material MAT141, following the same
# pattern as the existing row.
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "MAT141"
$ws.Range("C3").Value = "store/materials/MAT141/conditional probability.pdf"
